# Applies odds updates to Sheet1 as described in the commit diff
# "Atualizando o arquivo XLSX"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("N6").Value = 12

# Row 11
$ws.Range("G11").Value = 2.2
$ws.Range("H11").Value = 3.3
$ws.Range("I11").Value = 3.2
$ws.Range("L11").Value = 3.75
$ws.Range("AI11").Value = 15
$ws.Range("AL11").Value = 26
$ws.Range("AN11").Value = 4.33
$ws.Range("AZ11").Value = 51

# Row 12
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 3.7
$ws.Range("I12").Value = 1.65
$ws.Range("J12").Value = 5.1
$ws.Range("K12").Value = 2.2
$ws.Range("L12").Value = 2.18
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 7.2
$ws.Range("O12").Value = 1.31
$ws.Range("P12").Value = 3.15
$ws.Range("Q12").Value = 1.95
$ws.Range("R12").Value = 1.8
$ws.Range("S12").Value = 1.38
$ws.Range("T12").Value = 2.8
$ws.Range("U12").Value = 1.93
$ws.Range("V12").Value = 1.78
$ws.Range("W12").Value = 12.5
$ws.Range("X12").Value = 28
$ws.Range("Y12").Value = 16
$ws.Range("Z12").Value = 90
$ws.Range("AC12").Value = 7.2
$ws.Range("AD12").Value = 7.1
$ws.Range("AE12").Value = 17.5
$ws.Range("AH12").Value = 6.3
$ws.Range("AI12").Value = 7.2
$ws.Range("AK12").Value = 12
$ws.Range("AL12").Value = 13.5
$ws.Range("AN12").Value = 6.6
$ws.Range("AO12").Value = 28
$ws.Range("AT12").Value = 2.8
$ws.Range("AU12").Value = 7.7
$ws.Range("AW12").Value = 3.45
$ws.Range("AX12").Value = 7.9
$ws.Range("AY12").Value = 18
$ws.Range("AZ12").Value = 26

# Row 13
$ws.Range("G13").Value = 6.25
$ws.Range("I13").Value = 1.38
$ws.Range("N13").Value = 15
$ws.Range("W13").Value = 21
$ws.Range("X13").Value = 41
$ws.Range("Y13").Value = 21
$ws.Range("Z13").Value = 81
$ws.Range("AX13").Value = 6.5

# Row 14
$ws.Range("H14").Value = 3.25
$ws.Range("M14").Value = 1.11
$ws.Range("N14").Value = 6.5
$ws.Range("U14").Value = 2.25
$ws.Range("V14").Value = 1.57
$ws.Range("AB14").Value = 41
$ws.Range("AC14").Value = 6.5
$ws.Range("AJ14").Value = 19
$ws.Range("AN14").Value = 3.5
$ws.Range("AZ14").Value = 126

# Row 15
$ws.Range("H15").Value = 3.7
$ws.Range("J15").Value = 2.38
$ws.Range("U15").Value = 1.8
$ws.Range("V15").Value = 1.91
$ws.Range("W15").Value = 7.5
$ws.Range("Z15").Value = 13
$ws.Range("AC15").Value = 11
$ws.Range("AH15").Value = 13
$ws.Range("AO15").Value = 9
$ws.Range("AQ15").Value = 29
$ws.Range("AU15").Value = 8
$ws.Range("AW15").Value = 6.5
$ws.Range("AY15").Value = 29
$ws.Range("BB15").Value = 201

# Row 16
$ws.Range("G16").Value = 2.25
$ws.Range("I16").Value = 3.2
$ws.Range("J16").Value = 3
$ws.Range("W16").Value = 8
$ws.Range("AI16").Value = 15
$ws.Range("AO16").Value = 13
$ws.Range("AX16").Value = 17

# Row 19
$ws.Range("G19").Value = 1.75
$ws.Range("I19").Value = 4.5
$ws.Range("N19").Value = 10
$ws.Range("Q19").Value = 2.08
$ws.Range("R19").Value = 1.73
$ws.Range("Z19").Value = 13
$ws.Range("AL19").Value = 41
$ws.Range("AW19").Value = 6.5
$ws.Range("BA19").Value = 126

